# Apply "Add drag & lift values to excel sheet" change
$wb = $excel.ActiveWorkbook

# Add the new worksheet after Sheet1
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Drag & Lift"

# Header row
$newSheet.Range("A1").Value = "Airspeed (m/s)"
$newSheet.Range("B1").Value = "Lift (N)"
$newSheet.Range("C1").Value = "Drag (N)"

# Data rows
$data = @(
    @(0, 0, 0),
    @(4, 4.8, 0.6),
    @(8, 20, 2.2999999999999998),
    @(12, 47.5, 5.15),
    @(16, 86.4, 9.1)
)

$row = 2
foreach ($r in $data) {
    $newSheet.Cells.Item($row, 1).Value = $r[0]
    $newSheet.Cells.Item($row, 2).Value = $r[1]
    $newSheet.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Column widths to match the target layout
$newSheet.Columns.Item(1).ColumnWidth = 15.1796875
$newSheet.Columns.Item(2).ColumnWidth = 12.90625
$newSheet.Columns.Item(3).ColumnWidth = 13.7265625
$newSheet.Columns.Item(4).ColumnWidth = 8.7265625
$newSheet.Columns.Item(5).ColumnWidth = 8.7265625

# Select B6 as the active cell on the new sheet, and make it the active tab
$newSheet.Range("B6").Select()
$newSheet.Activate()
